# Applies the "Ohtani" error-handling fix: swaps the three featured stacks
# on Sheet1 (Seattle Mariners hitters / Atlanta Braves hitters / Colorado
# Rockies righties) and clears out the now-defunct fourth (Chicago White Sox)
# stack's player rows so the div/0 guard can do its job.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---- Table 1 (columns A:D) - Seattle Mariners hitters (FD, DK) ----
$ws1.Range("A1").Value = "Seattle Mariners hitters (FD, DK)"

$ws1.Range("A3").Value = "Segura"
$ws1.Range("B3").Value = 3300
$ws1.Range("C3").Value = 15.7

$ws1.Range("A4").Value = "Heredia"
$ws1.Range("B4").Value = 2200
$ws1.Range("C4").Value = 0

$ws1.Range("A5").Value = "Haniger"
$ws1.Range("B5").Value = 3500
$ws1.Range("C5").Value = 9.5

$ws1.Range("A6").Value = "Healy"
$ws1.Range("B6").Value = 2800
$ws1.Range("C6").Value = 0
$ws1.Range("D6").Formula = "=(C6 / B6) * 1000"

# ---- Table 2 (columns F:I) - Atlanta Braves hitters (FD, DK) ----
$ws1.Range("F1").Value = "Atlanta Braves hitters (FD, DK)"

$ws1.Range("F3").Value = "Inciarte"
$ws1.Range("G3").Value = 2800
$ws1.Range("H3").Value = 0

$ws1.Range("F4").Value = "Freeman"
$ws1.Range("G4").Value = 4300
$ws1.Range("H4").Value = 6

$ws1.Range("F5").Value = "Markakis"
$ws1.Range("G5").Value = 3400
$ws1.Range("H5").Value = 0

$ws1.Range("F6").Value = "Suzuki"
$ws1.Range("G6").Value = 2500
$ws1.Range("H6").Value = 0
$ws1.Range("I6").Formula = "=(H6 / G6) * 1000"

# ---- Table 3 (columns K:N) - Colorado Rockies righties (FD, DK) ----
$ws1.Range("K1").Value = "Colorado Rockies righties (FD, DK)"

$ws1.Range("K3").Value = "Cuevas"
$ws1.Range("L3").Value = 3000
$ws1.Range("M3").Value = 6.2

$ws1.Range("K4").Value = "Arenado"
$ws1.Range("L4").Value = 5400
$ws1.Range("M4").Value = 28.7

$ws1.Range("K5").Value = "Story"
$ws1.Range("L5").Value = 4400
$ws1.Range("M5").Value = 6

$ws1.Range("K6").Value = "Desmond"
$ws1.Range("L6").Value = 3700
$ws1.Range("M6").Value = 6

# Row 8 labels flip from Success to Failure for all three stacks now that
# none of them cleared value.
$ws1.Range("D8").Value = "Failure"
$ws1.Range("I8").Value = "Failure"
$ws1.Range("N8").Value = "Failure"

# ---- Fourth stack (Chicago White Sox hitters) retired; clear its header
# and player rows so the summary table falls back to the #DIV/0! guard
# (handles pitchers like Ohtani who also start as hitters). ----
$ws1.Range("A10").Value = ""

$ws1.Range("A12").Value = ""
$ws1.Range("B12").Value = ""
$ws1.Range("C12").Value = ""

$ws1.Range("A13").Value = ""
$ws1.Range("B13").Value = ""
$ws1.Range("C13").Value = ""

$ws1.Range("A14").Value = ""
$ws1.Range("B14").Value = ""
$ws1.Range("C14").Value = ""

$ws1.Range("A15").Value = ""
$ws1.Range("B15").Value = ""
$ws1.Range("C15").Value = ""

# Restore the F16/G16/H16 "Total:" row (mirrors the other three tables) so
# the now-empty fourth stack also reports 0s instead of staying blank.
$ws1.Range("F16").Value = "Total:"
$ws1.Range("G16").Formula = "=SUM(G12:G15)"
$ws1.Range("H16").Formula = "=SUM(H12:H15)"

# D17's Success/Failure label is cleared along with the rest of the fourth
# stack (its Total row is now a #DIV/0!).
$ws1.Range("D17").Value = ""

# E21 footnote text is unchanged in content but shifts shared-string index;
# re-assert it explicitly for safety.
$ws1.Range("E21").Value = "    "

# ---- Selections ---- (select sheet2 first so the final active tab/selection
# lands back on Sheet1, matching tabSelected="1" staying on Sheet1)
$ws2.Range("A1:A3").Select()
$ws1.Range("C6").Select()

$wb.Save()
